$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9121159911155701
$ws.Range("B1").Value = 2.725085496902466
$ws.Range("C1").Value = 4.447235107421875
$ws.Range("D1").Value = 2.155037879943848
$ws.Range("E1").Value = 1.273527264595032
